$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163

# --- Step 1: capture the values we will need to re-place after the row
#     shift below, while the original row numbers are still in effect. ---
$teacherName = $ws.Range("B13").Value()      # "5701460 - Antonio Iacono"

$metodoB = $ws.Range("B19").Value()          # "Provas e Trabalhos"
$metodoC = $ws.Range("C19").Value()
$criterioB = $ws.Range("B20").Value()        # "M = (0,6P + 0,4T)..."
$criterioC = $ws.Range("C20").Value()
$normaB = $ws.Range("B21").Value()           # "MF = (0,5 M + 0,5 R)..."
$normaC = $ws.Range("C21").Value()

# --- Step 2: remove the old "Docentes responsaveis" value row (row 13).
#     This shifts rows 14-24 up into 13-23 (labels in column A, and row
#     heights, move with it); column B/C content for several rows must
#     then be explicitly restored/relocated below since it does not track
#     that same shift semantically. ---
$ws.Rows(13).Delete()

# --- Step 3: write the new / relocated values into their final positions ---

# Objetivos: value becomes the teacher name
$ws.Range("B10").Value = $teacherName
$ws.Range("C10").Value = $teacherName

# Programa resumido: value becomes "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Programa: value becomes the activation date. Copy/PasteSpecial(values)
# keeps it stored as literal text instead of Excel re-parsing the
# date-shaped string into a date serial number.
$ws.Range("B8:C8").Copy()
$ws.Range("B15:C15").PasteSpecial($xlPasteValues)
$excel.CutCopyMode = 0

# Metodo: value becomes the teacher name again
$ws.Range("B18").Value = $teacherName
$ws.Range("C18").Value = $teacherName

# Critério / Norma de recuperação / Bibliografia rows keep the text that
# used to live at that same row number (only the column-A label shifted
# up), so restore it from what we captured in Step 1.
$ws.Range("B19").Value = $metodoB
$ws.Range("C19").Value = $metodoC
$ws.Range("B20").Value = $criterioB
$ws.Range("C20").Value = $criterioC
$ws.Range("B21").Value = $normaB
$ws.Range("C21").Value = $normaC
